# Updated symbol list on Mon Dec 26 05:59:02 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> new text value. All these cells are stored as inline strings
# (t="inlineStr") in the source workbook, so we force text assignment
# (NumberFormat "@") to avoid Excel silently re-casting the text as a number
# and stripping/altering formatting (e.g. trailing zeros like "0.7801" or
# "0.003000").

$updates = @{
    "D2"  = "243.75"
    "D3"  = "23.05"
    "D4"  = "5.417"
    "D5"  = "0.05959"
    "D6"  = "3.454"
    "D7"  = "6.517"
    "D8"  = "0.8134"
    "D9"  = "0.9203"
    "D10" = "0.1409"
    "D11" = "0.07443"
    "D12" = "0.03271"
    "D13" = "0.03059"
    "D14" = "0.09354"
    "D15" = "3.847"
    "D16" = "0.001568"
    "D17" = "0.04673"
    "D18" = "0.0005941"
    "D19" = "0.006078"
    "D20" = "0.004998"
    "E20" = "19HotbitTokenHTB"
    "D21" = "0.0009814"
    "D22" = "0.00007801"
    "D23" = "3.606"
    "D24" = "2.150"
    "D27" = "0.0002395"
    "D40" = "0.03943"
    "D41" = "0.006174"
    "D42" = "0.1073"
    "D43" = "0.003000"
    "D44" = "0.008449"
    "E44" = "43LocalTradersLCTBestin24h"
    "D45" = "0.00005241"
    "D48" = "0.7801"
    "D49" = "0.002289"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Temporarily force text format so Excel does not reinterpret the
    # numeric-looking string and strip/alter trailing zeros (e.g.
    # "0.7800" -> "0.78" or "0.003000" -> "3E-3"). Afterwards restore the
    # cell's style to "Normal" so no residual number-format / style index
    # is left behind (the source file has no explicit style on these cells).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
